$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate cells in the same order the strings were first authored so the
# resulting shared-strings table lines up with the source workbook.
$ws.Range("B1").Value = "Topic"
$ws.Range("C1").Value = "Problem Desc"
$ws.Range("B2").Value = "Heap"
$ws.Range("C2").Value = "Heap ADT"
$ws.Range("D1").Value = "Problem Link/Details"
$ws.Range("E1").Value = "ClassName"
$ws.Range("E2").Value = "HeapImpl.java"
$ws.Range("A1").Value = "SerialNo"
$ws.Range("A2").Value = 1

# Column widths (closest achievable character widths)
$ws.Columns.Item(2).ColumnWidth = 16
$ws.Columns.Item(3).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 21.5

# Selection / active cell
$ws.Range("A3").Select() | Out-Null
